# tests/data/validate_tables/valid.xlsx -- "added cell type checks"
#
# On the "params" sheet, row 3 (the "b" / "interp" row) gains two new
# cells describing an interpolation type and its value map, a new
# column (D) is widened to fit that value, and the active selection
# moves to E8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# New cells on row 3: C3 = "linear", D3 = the value JSON.
$ws.Range("C3").Value = "linear"
$ws.Range("D3").Value = '{"2020-01-01":10, "2031-06-01":9.5}'

# Match the formatting already used on the other header-ish cells
# (F1, I1, J1, K1, L1, M1) for the new C3 cell.
$ws.Range("C3").Font.Name = $ws.Range("F1").Font.Name
$ws.Range("C3").Font.Size = $ws.Range("F1").Font.Size
$ws.Range("C3").Font.Bold = $ws.Range("F1").Font.Bold
$ws.Range("C3").Font.Italic = $ws.Range("F1").Font.Italic
$ws.Range("C3").Font.Color = $ws.Range("F1").Font.Color

# Widen the new column D so the long JSON value is readable.
$ws.Columns.Item(4).ColumnWidth = 29.6666666

# Move the selection, as recorded in the saved view state.
$ws.Range("E8").Select()
